# Fixed some errors in Aerodynamic and Stability
# Updates computed balance values on the "GLOBAL RESULTS" and
# "LANDING GEARS" sheets of the Balance.xlsx workbook.

$wb = $excel.ActiveWorkbook

# ---- GLOBAL RESULTS sheet (column C values) ----
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")

$wsGlobal.Range("C2").Value  = 59.73179055824007
$wsGlobal.Range("C3").Value  = 12.558592054199604
$wsGlobal.Range("C4").Value  = 23.881283765287307
$wsGlobal.Range("C5").Value  = 0.5489100388884587

$wsGlobal.Range("C7").Value  = 43.11606869451627
$wsGlobal.Range("C8").Value  = 12.176680569707568
$wsGlobal.Range("C9").Value  = 26.608857770674003
$wsGlobal.Range("C10").Value = 0.6116031825269179

$wsGlobal.Range("C12").Value = 43.11606869451627
$wsGlobal.Range("C13").Value = 12.176680569707568
$wsGlobal.Range("C14").Value = 26.608857770674003
$wsGlobal.Range("C15").Value = 0.6116031825269179

$wsGlobal.Range("C17").Value = 53.60136223687021
$wsGlobal.Range("C18").Value = 12.417684484851204
$wsGlobal.Range("C19").Value = 17.425612078158753
$wsGlobal.Range("C20").Value = 0.4005267680534276

$wsGlobal.Range("C22").Value = 51.624242905224534
$wsGlobal.Range("C23").Value = 12.372240502028212
$wsGlobal.Range("C24").Value = 24.64235345179721
$wsGlobal.Range("C25").Value = 0.5664031852085996

$wsGlobal.Range("C27").Value = 0.29269052532561396
$wsGlobal.Range("C28").Value = 0.6551257388670622

# ---- LANDING GEARS sheet (column C value) ----
$wsLandingGears = $wb.Worksheets.Item("LANDING GEARS")

$wsLandingGears.Range("C2").Value = 12.32108109437063
